# Task 1 / subtask 2 — Machine Learning Course At Universities.xlsx
# Adds the per-university "total courses" column (M) as a SUM formula,
# records a new course link (K50), wires up two missing hyperlinks
# (I50, J50), and moves the saved view/selection to the new work area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- M27:M47 -> shared formula =SUM(C:L) for that row ------------------
$ws.Range("M27:M47").Formula = "=SUM(C27:L27)"

# --- M48 -> same sum, but outside the shared range above (own formula) -
$ws.Range("M48").Formula = "=SUM(C48:L48)"

# --- New course link for Columbia (column K) on row 50 -----------------
$ws.Range("K50").Value = "http://www.cs.columbia.edu/~jebara/4771/notes/class1x.pdf"

# --- Hyperlinks that were missing for the Michigan State / Notre Dame --
# --- syllabus links already typed as text in I50 / J50 -----------------
# (re-apply the original cell style afterwards -- Hyperlinks.Add mints a
# fresh "visited/hyperlink" style on the target cell, but these cells
# already carried that same look, so restore the untouched style id.)
$i50Style = $ws.Range("I50").Style
$j50Style = $ws.Range("J50").Style
$ws.Hyperlinks.Add($ws.Range("I50"), "https://github.com/jiayuzhou/CSE491-2016Fall")
$ws.Range("I50").Style = $i50Style
$ws.Hyperlinks.Add($ws.Range("J50"), "https://github.com/cse40625/cse40625/blob/master/syllabus.pdf")
$ws.Range("J50").Style = $j50Style

# --- Move the saved viewport / selection down to the new work area -----
$ws.Range("K52").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
